$p = $ppt.ActivePresentation

# --- Slide 2: fix punctuation "so data serves you" -> ", so data serves you" ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(15)
$tr2 = $sh2.TextFrame.TextRange
$run2 = $tr2.Runs(3)
$run2.Text = " today. But you will learn to ask the right questions, so data serves you, not the other way around."

# --- Slide 28: trim bullet text down to "Select one  cases " ---
$s28 = $p.Slides.Item(28)
$sh28 = $s28.Shapes.Item(6)
$sh28.TextFrame.TextRange.Text = "• Select one  cases "

# --- Slide 29: drop "for YEARS" from the key-insight callout ---
$s29 = $p.Slides.Item(29)
$sh29 = $s29.Shapes.Item(26)
$sh29.TextFrame.TextRange.Text = "Key Insight from Vulnerability Case: The Karachi Malir health crisis (2.4M people, 82nd percentile health) went unnoticed  because no one questioned aggregate metrics (41st overall). Your job: Ask the right questions BEFORE crises become visible."

# --- Slide 30: remove the stray empty teal bar shape ("Shape 18") ---
$s30 = $p.Slides.Item(30)
$sh30 = $s30.Shapes.Item(19)
$sh30.Delete()
